$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the two HU title strings (HU-07 and HU-08 rows)
$ws.Range("C14").Value = "Generación de carnet de vacunación"
$ws.Range("C15").Value = "Registro de información clínica-vacunas"

# Adjust row heights for rows 11-17
$ws.Rows.Item(11).RowHeight = 20.25
$ws.Rows.Item(12).RowHeight = 20.25
$ws.Rows.Item(13).RowHeight = 20.25
$ws.Rows.Item(14).RowHeight = 20.25
$ws.Rows.Item(15).RowHeight = 20.25
$ws.Rows.Item(16).RowHeight = 20.25
$ws.Rows.Item(17).RowHeight = 21
